$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'290.44"
$ws.Range("E2").Value = "'-3.44%"
$ws.Range("D3").Value = "'30.60"
$ws.Range("E3").Value = "'-6.04%"
$ws.Range("D4").Value = "'4.915"
$ws.Range("E4").Value = "'-3.05%"
$ws.Range("D5").Value = "'0.07246"
$ws.Range("E5").Value = "'-6.03%"
$ws.Range("D6").Value = "'1.810"
$ws.Range("E6").Value = "'-10.46%"
$ws.Range("D7").Value = "'7.640"
$ws.Range("E7").Value = "'-3.29%"
$ws.Range("D8").Value = "'3.697"
$ws.Range("E8").Value = "'-2.83%"
$ws.Range("D9").Value = "'0.9003"
$ws.Range("E9").Value = "'-2.78%"
$ws.Range("D10").Value = "'0.1683"
$ws.Range("E10").Value = "'-4.09%"
$ws.Range("D11").Value = "'0.08023"
$ws.Range("E11").Value = "'-2.10%"
$ws.Range("D12").Value = "'0.08055"
$ws.Range("E12").Value = "'-6.23%"
$ws.Range("D13").Value = "'0.03063"
$ws.Range("E13").Value = "'-1.28%"
$ws.Range("D14").Value = "'0.1002"
$ws.Range("E14").Value = "'0.22%"
$ws.Range("E15").Value = "'-1.57%"
$ws.Range("D16").Value = "'0.005749"
$ws.Range("E16").Value = "'-2.63%"
$ws.Range("D17").Value = "'3.476"
$ws.Range("E17").Value = "'-0.02%"
$ws.Range("E18").Value = "'-3.63%"
$ws.Range("D19").Value = "'0.3318"
$ws.Range("E19").Value = "'-0.49%"
$ws.Range("D20").Value = "'0.1304"
$ws.Range("E20").Value = "'-1.87%"
$ws.Range("D21").Value = "'3.955"
$ws.Range("E21").Value = "'-10.03%"
$ws.Range("E22").Value = "'9.52%"
$ws.Range("D23").Value = "'0.04505"
$ws.Range("D24").Value = "'0.001213"
$ws.Range("E24").Value = "'-1.43%"
$ws.Range("D25").Value = "'0.004427"
$ws.Range("E25").Value = "'7.15%"
$ws.Range("D26").Value = "'0.0001298"
$ws.Range("E26").Value = "'3.49%"
$ws.Range("D27").Value = "'0.0003394"
$ws.Range("E27").Value = "'-95.47%"
$ws.Range("D39").Value = "'0.01589"
$ws.Range("E39").Value = "'-7.88%"
$ws.Range("D40").Value = "'0.04348"
$ws.Range("E40").Value = "'-7.25%"
$ws.Range("D41").Value = "'0.007322"
$ws.Range("E41").Value = "'-2.45%"
$ws.Range("D42").Value = "'0.01004"
$ws.Range("D43").Value = "'0.1313"
$ws.Range("E43").Value = "'-3.31%"
$ws.Range("D44").Value = "'0.002001"
$ws.Range("E44").Value = "'-10.70%"
$ws.Range("D45").Value = "'0.009464"
$ws.Range("E45").Value = "'-10.20%"
$ws.Range("D46").Value = "'0.00005852"
$ws.Range("E46").Value = "'-5.07%"
$ws.Range("D47").Value = "'0.00000000749"
$ws.Range("E47").Value = "'-0.47%"
$ws.Range("D48").Value = "'2.255"
$ws.Range("E48").Value = "'43.53%"
$ws.Range("D49").Value = "'0.002898"
$ws.Range("E49").Value = "'17.94%"
$ws.Range("E50").Value = "'-0.47%"
$ws.Range("E51").Value = "'-0.47%"
